$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7085.032
$ws.Range("I19").Value = 789.8570999999999
$ws.Range("J19").Value = 8921.125
$ws.Range("K19").Value = 789.8570999999999
$ws.Range("L19").Value = 8921.125
$ws.Range("M19").Value = -614.8570999999999
$ws.Range("N19").Value = -9271.125
$ws.Range("H107").Value = 1358.7693
$ws.Range("I107").Value = 1505.8182
$ws.Range("J107").Value = 550
$ws.Range("K107").Value = 1505.8182
$ws.Range("L107").Value = 550
$ws.Range("M107").Value = 414.1818000000001
$ws.Range("N107").Value = -4390
$ws.Range("H137").Value = 3487.3823
$ws.Range("I137").Value = 2419.3667
$ws.Range("J137").Value = 11497.5
$ws.Range("K137").Value = 7258.1001
$ws.Range("L137").Value = 34492.5
$ws.Range("M137").Value = -4708.1001
$ws.Range("N137").Value = -39592.5
$ws.Range("H138").Value = 2537.4658
$ws.Range("J138").Value = 2470.2703
$ws.Range("L138").Value = 7410.8109
$ws.Range("N138").Value = -17690.8109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3142.5173
$ws.Range("I61").Value = 2711.9375
$ws.Range("J61").Value = 3672.4614
$ws.Range("K61").Value = 2711.9375
$ws.Range("L61").Value = 3672.4614
$ws.Range("M61").Value = -2499.9375
$ws.Range("N61").Value = -4096.4614
$ws.Range("H74").Value = 2585
$ws.Range("J74").Value = 3666
$ws.Range("L74").Value = 3666
$ws.Range("N74").Value = -5414
$ws.Range("H77").Value = 2585
$ws.Range("J77").Value = 3666
$ws.Range("L77").Value = 18330
$ws.Range("N77").Value = -27066
$ws.Range("H92").Value = 71000
$ws.Range("J92").Value = 71000
$ws.Range("L92").Value = 71000
$ws.Range("N92").Value = -75992
$ws.Range("H101").Value = 74512.44500000001
$ws.Range("J101").Value = 74512.44500000001
$ws.Range("L101").Value = 74512.44500000001
$ws.Range("N101").Value = -81002.44500000001
$ws.Range("H132").Value = 5201.6333
$ws.Range("I132").Value = 4401
$ws.Range("K132").Value = 13203
$ws.Range("M132").Value = -10673
$ws.Range("H136").Value = 3142.5173
$ws.Range("I136").Value = 2711.9375
$ws.Range("J136").Value = 3672.4614
$ws.Range("K136").Value = 8135.8125
$ws.Range("L136").Value = 11017.3842
$ws.Range("M136").Value = -5585.8125
$ws.Range("N136").Value = -16117.3842
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 100000
$ws.Range("J100").Value = 100000
$ws.Range("L100").Value = 100000
$ws.Range("N100").Value = -102164
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5965.1064
$ws.Range("I31").Value = 1214.6666
$ws.Range("J31").Value = 12378.2
$ws.Range("K31").Value = 1214.6666
$ws.Range("L31").Value = 12378.2
$ws.Range("M31").Value = -919.6666
$ws.Range("N31").Value = -12968.2
$ws.Range("H34").Value = 5965.1064
$ws.Range("I34").Value = 1214.6666
$ws.Range("J34").Value = 12378.2
$ws.Range("K34").Value = 1214.6666
$ws.Range("L34").Value = 12378.2
$ws.Range("M34").Value = -1012.6666
$ws.Range("N34").Value = -12782.2
$ws.Range("H58").Value = 1580.5416
$ws.Range("I58").Value = 1440.5555
$ws.Range("K58").Value = 1440.5555
$ws.Range("M58").Value = -1237.5555
$ws.Range("H74").Value = 25745.25
$ws.Range("J74").Value = 25745.25
$ws.Range("L74").Value = 25745.25
$ws.Range("N74").Value = -27493.25
$ws.Range("H77").Value = 25745.25
$ws.Range("J77").Value = 25745.25
$ws.Range("L77").Value = 77235.75
$ws.Range("N77").Value = -85971.75
$ws.Range("H88").Value = 34333.332
$ws.Range("J88").Value = 34333.332
$ws.Range("L88").Value = 34333.332
$ws.Range("N88").Value = -35145.332
$ws.Range("H91").Value = 34333.332
$ws.Range("J91").Value = 34333.332
$ws.Range("L91").Value = 34333.332
$ws.Range("N91").Value = -37141.332
$ws.Range("H96").Value = 42450
$ws.Range("J96").Value = 42450
$ws.Range("L96").Value = 42450
$ws.Range("N96").Value = -47942
$ws.Range("H132").Value = 8774624
$ws.Range("I132").Value = 2548
$ws.Range("K132").Value = 7644
$ws.Range("M132").Value = -5114
$ws.Range("H134").Value = 1917.6666
$ws.Range("I134").Value = 1752.7333
$ws.Range("J134").Value = 2330
$ws.Range("K134").Value = 5258.199900000001
$ws.Range("L134").Value = 6990
$ws.Range("M134").Value = -2723.199900000001
$ws.Range("N134").Value = -12060
$ws.Range("H136").Value = 1580.5416
$ws.Range("I136").Value = 1440.5555
$ws.Range("K136").Value = 4321.666499999999
$ws.Range("M136").Value = -1771.666499999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 323.66666
$ws.Range("I5").Value = 323.66666
$ws.Range("K5").Value = 970.9999799999999
$ws.Range("M5").Value = -858.9999799999999
$ws.Range("H113").Value = 990.4231
$ws.Range("I113").Value = 572.5
$ws.Range("J113").Value = 1066.409
$ws.Range("K113").Value = 1717.5
$ws.Range("L113").Value = 3199.227
$ws.Range("M113").Value = 452.5
$ws.Range("N113").Value = -7539.227000000001
$ws.Range("H122").Value = 6797.6875
$ws.Range("I122").Value = 426.5
$ws.Range("J122").Value = 17416.334
$ws.Range("K122").Value = 3838.5
$ws.Range("L122").Value = 156747.006
$ws.Range("M122").Value = -1388.5
$ws.Range("N122").Value = -161647.006
$ws.Range("H135").Value = 323.66666
$ws.Range("I135").Value = 323.66666
$ws.Range("K135").Value = 2912.99994
$ws.Range("M135").Value = -377.9999399999997
$ws.Range("H140").Value = 2051.9333
$ws.Range("I140").Value = 1815.75
$ws.Range("J140").Value = 2321.8572
$ws.Range("K140").Value = 5447.25
$ws.Range("L140").Value = 6965.571599999999
$ws.Range("M140").Value = -267.25
$ws.Range("N140").Value = -17325.5716
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 143769.72
$ws.Range("I113").Value = 250472
$ws.Range("K113").Value = 250472
$ws.Range("M113").Value = -248302
$ws.Range("H132").Value = 3126.7368
$ws.Range("I132").Value = 3298.8572
$ws.Range("J132").Value = 3026.3333
$ws.Range("K132").Value = 9896.571599999999
$ws.Range("L132").Value = 9078.999899999999
$ws.Range("M132").Value = -7366.571599999999
$ws.Range("N132").Value = -14138.9999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 31030
$ws.Range("J104").Value = 31030
$ws.Range("L104").Value = 31030
$ws.Range("N104").Value = -38018
$ws.Range("H132").Value = 2567.4688
$ws.Range("I132").Value = 1550.6316
$ws.Range("K132").Value = 4651.8948
$ws.Range("M132").Value = -2121.8948
$ws.Range("H136").Value = 11906774
$ws.Range("I136").Value = 1717.1111
$ws.Range("K136").Value = 5151.3333
$ws.Range("M136").Value = -2601.3333
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 11092.25
$ws.Range("I45").Value = 7569
$ws.Range("J45").Value = 12266.667
$ws.Range("K45").Value = 7569
$ws.Range("L45").Value = 12266.667
$ws.Range("M45").Value = -7078
$ws.Range("N45").Value = -13248.667
$ws.Range("H101").Value = 18700.666
$ws.Range("J101").Value = 18700.666
$ws.Range("L101").Value = 18700.666
$ws.Range("N101").Value = -25190.666
$ws.Range("H104").Value = 47666.668
$ws.Range("J104").Value = 47666.668
$ws.Range("L104").Value = 47666.668
$ws.Range("N104").Value = -54654.668
$ws.Range("H132").Value = 3970465
$ws.Range("I132").Value = 2242.8262
$ws.Range("J132").Value = 8774103
$ws.Range("K132").Value = 6728.4786
$ws.Range("L132").Value = 26322309
$ws.Range("M132").Value = -4198.4786
$ws.Range("N132").Value = -26327369
$ws.Range("H136").Value = 2544.5405
$ws.Range("I136").Value = 2523.1738
$ws.Range("K136").Value = 7569.5214
$ws.Range("M136").Value = -5019.5214
